$d = $word.ActiveDocument

# --- 1. Rewrite the ranking list -------------------------------------------
# Paragraph 1 (Title) stays untouched.
# Paragraph 2 (was the 🥇 "Toast" Heading1 line) becomes the centered Subtitle
# caption. Paragraphs 3-16 keep their relative order but get new text,
# shift to the "Normal" style, and are centered - matching the target diff.
$items = @(
    @{ Index = 2;  Text = "(Dish price / minute of preparation)";    Style = "Subtitle" },
    @{ Index = 3;  Text = "🥇 Steak profit: 7.35 ₪/min";             Style = "Normal"   },
    @{ Index = 4;  Text = "🥈 Stuffed Mushrooms profit: 6.81 ₪/min"; Style = "Normal"   },
    @{ Index = 5;  Text = "🥉 Salad profit: 6.00 ₪/min";             Style = "Normal"   },
    @{ Index = 6;  Text = "4) Brule Cream profit: 5.36 ₪/min";       Style = "Normal"   },
    @{ Index = 7;  Text = "5) Pasta profit: 4.65 ₪/min";             Style = "Normal"   },
    @{ Index = 8;  Text = "6) Pizza profit: 3.63 ₪/min";             Style = "Normal"   },
    @{ Index = 9;  Text = "7) Krep profit: 3.56 ₪/min";              Style = "Normal"   },
    @{ Index = 10; Text = "8) Belgian Waffle profit: 3.20 ₪/min";    Style = "Normal"   },
    @{ Index = 11; Text = "9) Hamburger profit: 2.70 ₪/min";         Style = "Normal"   },
    @{ Index = 12; Text = "10) Empanadas profit: 2.68 ₪/min";        Style = "Normal"   },
    @{ Index = 13; Text = "11) Schnitzel profit: 2.48 ₪/min";        Style = "Normal"   },
    @{ Index = 14; Text = "12) Cake profit: 2.32 ₪/min";             Style = "Normal"   },
    @{ Index = 15; Text = "13) Roast profit: 2.00 ₪/min";            Style = "Normal"   },
    @{ Index = 16; Text = "14) Arancini profit: 1.81 ₪/min";         Style = "Normal"   }
)

foreach ($item in $items) {
    $p = $d.Paragraphs.Item($item.Index)
    $p.Range.Text = $item.Text
    $p.Style = $item.Style
    $p.Alignment = 1   # wdAlignParagraphCenter
}

# --- 2. Make sure the "List Paragraph" style exists in the style catalogue -
# The target stylesheet gains a (built-in, currently unused) "List Paragraph"
# style definition. Applying the built-in style to a throwaway paragraph and
# then removing that paragraph causes Word to "mint" the style definition
# into the document's style catalogue without leaving any visible content
# behind, exactly mirroring how Word itself adds latent built-in styles the
# first time they are touched.
$d.Paragraphs.Add() | Out-Null
$scratch = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratch.Range.Text = "scratch"
$scratch.Style = "List Paragraph"
$scratch.Range.Delete()

$listStyle = $d.Styles.Item("List Paragraph")
$listStyle.Priority = 34
$listStyle.NoSpaceBetweenParagraphsOfSameStyle = $true
$listStyle.ParagraphFormat.LeftIndent = 36
